# Scheduled-runner refresh of market/profit figures across the Leve sheets.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) for a
# set of rows on each class sheet, including a couple of rows whose NQ/HQ
# profit cells go from blank to populated (or vice versa) as prices move to/
# from 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1006
$ws.Range("J125").Value = 1359
$ws.Range("L125").Value = 12231
$ws.Range("N125").Value = -17151

$ws.Range("H137").Value = 2633439.8
$ws.Range("I137").Value = 4168341.2
$ws.Range("K137").Value = 12505023.6
$ws.Range("M137").Value = -12502473.6

$ws.Range("H138").Value = 1955461.4
$ws.Range("I138").Value = 114121.555
$ws.Range("J138").Value = 2170683
$ws.Range("K138").Value = 342364.665
$ws.Range("L138").Value = 6512049
$ws.Range("M138").Value = -337224.665
$ws.Range("N138").Value = -6522329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7602494.5
$ws.Range("I32").Value = 9837220
$ws.Range("J32").Value = 4427.8
$ws.Range("K32").Value = 9837220
$ws.Range("L32").Value = 4427.8
$ws.Range("M32").Value = -9836933
$ws.Range("N32").Value = -5001.8

$ws.Range("H61").Value = 23304322
$ws.Range("I61").Value = 31282932
$ws.Range("J61").Value = 93820.73
$ws.Range("K61").Value = 31282932
$ws.Range("L61").Value = 93820.73
$ws.Range("M61").Value = -31282720
$ws.Range("N61").Value = -94244.73

$ws.Range("H74").Value = 7412748
$ws.Range("J74").Value = 78538.46000000001
$ws.Range("L74").Value = 78538.46000000001
$ws.Range("N74").Value = -80286.46000000001

$ws.Range("H77").Value = 7412748
$ws.Range("J77").Value = 78538.46000000001
$ws.Range("L77").Value = 392692.3
$ws.Range("N77").Value = -401428.3

$ws.Range("H110").Value = 3334333.2
$ws.Range("I110").Value = 3334333.2
$ws.Range("K110").Value = 3334333.2
$ws.Range("M110").Value = -3332288.2

$ws.Range("H132").Value = 60551.176
$ws.Range("I132").Value = 41581.68
$ws.Range("J132").Value = 113244.22
$ws.Range("K132").Value = 124745.04
$ws.Range("L132").Value = 339732.66
$ws.Range("M132").Value = -122215.04
$ws.Range("N132").Value = -344792.66

$ws.Range("H136").Value = 23304322
$ws.Range("I136").Value = 31282932
$ws.Range("J136").Value = 93820.73
$ws.Range("K136").Value = 93848796
$ws.Range("L136").Value = 281462.19
$ws.Range("M136").Value = -93846246
$ws.Range("N136").Value = -286562.19

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 68000
$ws.Range("J59").Value = 68000
$ws.Range("L59").Value = 68000
$ws.Range("N59").Value = -69694

$ws.Range("H107").Value = 3157.3914
$ws.Range("I107").Value = 3003.3076
$ws.Range("J107").Value = 3357.7
$ws.Range("K107").Value = 3003.3076
$ws.Range("L107").Value = 3357.7
$ws.Range("M107").Value = -1083.3076
$ws.Range("N107").Value = -7197.7

$ws.Range("H134").Value = 4238.1665
$ws.Range("I134").Value = 3005.087
$ws.Range("K134").Value = 9015.261
$ws.Range("M134").Value = -6480.261

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8193.415999999999
$ws.Range("I31").Value = 30043.37
$ws.Range("K31").Value = 30043.37
$ws.Range("M31").Value = -29748.37

$ws.Range("H34").Value = 8193.415999999999
$ws.Range("I34").Value = 30043.37
$ws.Range("K34").Value = 30043.37
$ws.Range("M34").Value = -29841.37

$ws.Range("H99").Value = 3331.5
$ws.Range("I99").Value = 4464.125
$ws.Range("J99").Value = 1821.3334
$ws.Range("K99").Value = 4464.125
$ws.Range("L99").Value = 1821.3334
$ws.Range("M99").Value = -2966.125
$ws.Range("N99").Value = -4817.3334

$ws.Range("H126").Value = 3331.5
$ws.Range("I126").Value = 4464.125
$ws.Range("J126").Value = 1821.3334
$ws.Range("K126").Value = 13392.375
$ws.Range("L126").Value = 5464.0002
$ws.Range("M126").Value = -10922.375
$ws.Range("N126").Value = -10404.0002

$ws.Range("H132").Value = 102880.4
$ws.Range("I132").Value = 1578
$ws.Range("J132").Value = 170415.33
$ws.Range("K132").Value = 4734
$ws.Range("L132").Value = 511245.99
$ws.Range("M132").Value = -2204
$ws.Range("N132").Value = -516305.99

$ws.Range("H134").Value = 32217.457
$ws.Range("I134").Value = 1299.5217
$ws.Range("J134").Value = 91476.836
$ws.Range("K134").Value = 3898.5651
$ws.Range("L134").Value = 274430.508
$ws.Range("M134").Value = -1363.5651
$ws.Range("N134").Value = -279500.508

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 821.76666
$ws.Range("I68").Value = 600.2857
$ws.Range("J68").Value = 889.1739
$ws.Range("K68").Value = 1800.8571
$ws.Range("L68").Value = 2667.5217
$ws.Range("M68").Value = -989.8571000000002
$ws.Range("N68").Value = -4289.5217

$ws.Range("H71").Value = 821.76666
$ws.Range("I71").Value = 600.2857
$ws.Range("J71").Value = 889.1739
$ws.Range("K71").Value = 5402.571300000001
$ws.Range("L71").Value = 8002.5651
$ws.Range("M71").Value = -1346.571300000001
$ws.Range("N71").Value = -16114.5651

$ws.Range("H131").Value = 867.2
$ws.Range("I131").Value = 531.1429000000001
$ws.Range("J131").Value = 969.4783
$ws.Range("K131").Value = 1593.4287
$ws.Range("L131").Value = 2908.4349
$ws.Range("M131").Value = 3446.5713
$ws.Range("N131").Value = -12988.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 339.875
$ws.Range("I2").Value = 428.16666
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 428.16666
$ws.Range("L2").Value = 75
$ws.Range("M2").Value = -315.16666
$ws.Range("N2").Value = -301

$ws.Range("H113").Value = 833.375
$ws.Range("I113").Value = 977.8333
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 977.8333
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = 1192.1667
$ws.Range("N113").Value = -4740

$ws.Range("H122").Value = 583.3333
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950

$ws.Range("H132").Value = 90230.78
$ws.Range("I132").Value = 64956.125
$ws.Range("J132").Value = 148001.42
$ws.Range("K132").Value = 194868.375
$ws.Range("L132").Value = 444004.26
$ws.Range("M132").Value = -192338.375
$ws.Range("N132").Value = -449064.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H132").Value = 33521.656
$ws.Range("I132").Value = 1745.55
$ws.Range("J132").Value = 86481.836
$ws.Range("K132").Value = 5236.65
$ws.Range("L132").Value = 259445.508
$ws.Range("M132").Value = -2706.65
$ws.Range("N132").Value = -264505.508

$ws.Range("H136").Value = 69006.44500000001
$ws.Range("I136").Value = 34928.9
$ws.Range("J136").Value = 147646.92
$ws.Range("K136").Value = 104786.7
$ws.Range("L136").Value = 442940.76
$ws.Range("M136").Value = -102236.7
$ws.Range("N136").Value = -448040.76

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 62704.727
$ws.Range("I132").Value = 47111.773
$ws.Range("J132").Value = 93890.63
$ws.Range("K132").Value = 141335.319
$ws.Range("L132").Value = 281671.89
$ws.Range("M132").Value = -138805.319
$ws.Range("N132").Value = -286731.89
